$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, shifting the existing rows 63-119 down to 64-120
$ws.Rows(63).Insert()

# Populate the new row 63 with the new weekly record (same dimension attrs as the
# record that is now in row 64, except for the date + volume/price/unit fields).
$ws.Range("A63").Value = 9
$ws.Range("B63").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44907
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = 100114002
$ws.Range("G63").Value = "Camote"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 790
$ws.Range("K63").Value = 17000
$ws.Range("L63").Value = 18000
$ws.Range("M63").Value = 17506
$ws.Range("N63").Value = "$/caja 18 kilos"
$ws.Range("O63").Value = "Perú"
$ws.Range("P63").Value = 973
$ws.Range("Q63").Value = 18
$ws.Range("R63").Value = "Hortaliza"
